$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet
$ws.Name = "Test Oligos Sheet 1"

# 2. Update the sequence value in C12 (KH11 / KHtest11 row)
$ws.Range("C12").Value = "AGAGGGGCTGGGAGTTGGACCCC"

# 3. Update the selected cell/range to B12
$null = $ws.Range("B12").Select()

# 4. Set the print scale on the page setup
$ws.PageSetup.Zoom = 62
